$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff swaps the Id (A), Antal (I), Ost (Q), and Nord (R) values
# between row 10 and row 11 while everything else in those rows stays
# the same. Use Copy/PasteSpecial (rather than re-typing .Value/.Value2)
# so cell types/formatting (e.g. the text-typed "Antal" column) and
# numeric precision are preserved exactly, with no stray style entries.

$scratch = $ws.Range("ZZ1")

# --- Column A (Id) ---
$ws.Range("A10").Copy()
$scratch.PasteSpecial()
$ws.Range("A11").Copy()
$ws.Range("A10").PasteSpecial()
$scratch.Copy()
$ws.Range("A11").PasteSpecial()

# --- Column I (Antal) ---
$ws.Range("I10").Copy()
$scratch.PasteSpecial()
$ws.Range("I11").Copy()
$ws.Range("I10").PasteSpecial()
$scratch.Copy()
$ws.Range("I11").PasteSpecial()

# --- Column Q (Ost) ---
$ws.Range("Q10").Copy()
$scratch.PasteSpecial()
$ws.Range("Q11").Copy()
$ws.Range("Q10").PasteSpecial()
$scratch.Copy()
$ws.Range("Q11").PasteSpecial()

# --- Column R (Nord) ---
$ws.Range("R10").Copy()
$scratch.PasteSpecial()
$ws.Range("R11").Copy()
$ws.Range("R10").PasteSpecial()
$scratch.Copy()
$ws.Range("R11").PasteSpecial()

$scratch.Clear()
